# Replace the "Latest Modified Debt Plot - PERSI" inline picture with a
# plain text run carrying the same wording, per the commit:
#   "Added modified version of the debtPlot (with title + grid lines)"
#
# The original run:
#   <w:r><w:drawing>...wp:docPr descr="Latest Modified Debt Plot - PERSI".../w:drawing></w:r>
# becomes:
#   <w:r><w:t xml:space="preserve">Latest Modified Debt Plot - PERSI</w:t></w:r>

$d = $word.ActiveDocument

$targetText = "Latest Modified Debt Plot - PERSI"

# Locate the inline picture by its alt text / description (the docPr descr
# attribute surfaces as InlineShape.AlternativeText), rather than assuming a
# fixed index, so the script is resilient to ordering.
$targetShape = $null
for ($i = 1; $i -le $d.InlineShapes.Count; $i++) {
    $shape = $d.InlineShapes.Item($i)
    if ($shape.AlternativeText -eq $targetText) {
        $targetShape = $shape
    }
}

if ($targetShape -ne $null) {
    # Remember where the picture's anchor character lives, then remove the
    # picture (and its owning run) entirely.
    $insertionPoint = $targetShape.Range.Start
    $targetShape.Delete()

    # Insert a fresh run with the same wording as plain text right where the
    # picture used to be (still inside the same paragraph, before the
    # trailing "Note that the..." paragraph).
    $ins = $d.Range($insertionPoint, $insertionPoint)
    $ins.InsertAfter($targetText)
}
